# Merge co-cooking variants (A/B/C) out of the relationship matrix.
# The three co_cooking_A / co_cooking_B / co_cooking_C rows+columns are
# removed entirely (both as a row and as a column) from the adjacency
# matrix on Sheet1; nothing else in the remaining grid changes value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns C, E, G hold co_cooking_A/B/C (header row 1).
# Rows 3, 5, 7 hold the same three entities (column A labels).
# Delete from the rightmost/bottommost first so earlier indices don't shift.

$ws.Range("G1:G23").EntireColumn.Delete() | Out-Null
$ws.Range("E1:E23").EntireColumn.Delete() | Out-Null
$ws.Range("C1:C23").EntireColumn.Delete() | Out-Null

$ws.Range("A7:A7").EntireRow.Delete() | Out-Null
$ws.Range("A5:A5").EntireRow.Delete() | Out-Null
$ws.Range("A3:A3").EntireRow.Delete() | Out-Null

# Restore the view state recorded after the edit.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("M24").Select()
